$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("O3").Value = 5
$ws.Range("R3").Value = 5
$ws.Range("V3").Value = 5

# Row 4
$ws.Range("O4").Value = 5
$ws.Range("R4").Value = 5
$ws.Range("U4").Value = 5
$ws.Range("V4").Value = 5

# Row 8
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 5
$ws.Range("U8").Value = 5
$ws.Range("V8").Value = 5

# Row 9
# N9 is a brand-new cell that needs the same (green-fill) style as R19,
# so copy its format first, then set the value.
$ws.Range("R19").Copy() | Out-Null
$ws.Range("N9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 5

# Row 11
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("U11").Value = 5
$ws.Range("V11").Value = 5

# Row 20
$ws.Range("R20").Value = 5
$ws.Range("T20").Value = 5
$ws.Range("U20").Value = 5

# Row 25
$ws.Range("N25").Value = 5
$ws.Range("O25").Value = 5
$ws.Range("P25").Value = 5
$ws.Range("U25").Value = 5
$ws.Range("V25").Value = 5

# Sheet view changes
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("N9").Select()
